# BodyPartDefs.xlsx — "adding Border.prefab, adding border to scene, and
# lowering some jump cooldowns"
#
# The spreadsheet-visible part of that commit is entirely on the
# ThoraxPartDefs sheet: every row's jumpCooldown (column J) drops to 1.
# The rest of the commit (Border.prefab / scene changes) lives outside this
# workbook. We also restore the per-sheet cell selections recorded in the
# saved view state for the three sheets that moved.

$wb = $excel.ActiveWorkbook

# --- Lower jump cooldowns on ThoraxPartDefs (column J = "jumpCooldown") ---
$thorax = $wb.Worksheets.Item("ThoraxPartDefs")
[void]$thorax.Activate()

for ($row = 2; $row -le 10; $row++) {
    $thorax.Cells.Item($row, 10).Value = 1
}

# --- Restore the remembered selection on each sheet ---
[void]$thorax.Range("H18").Select()

$abdomen = $wb.Worksheets.Item("AbdomenPartDefs")
[void]$abdomen.Activate()
[void]$abdomen.Range("K2").Select()

$leg = $wb.Worksheets.Item("LegPartDefs")
[void]$leg.Activate()
[void]$leg.Range("C4").Select()

# ThoraxPartDefs remains the active/visible tab.
[void]$thorax.Activate()
